# Season record columns (Wins / Losses / Ties) were missing from the
# exported team-statistics sheet; add them back.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells should look like the existing header row (bold font,
# thin border, centered/top alignment) - copy formatting from the last
# existing header cell (AC1) rather than re-declaring it by hand.
$ws.Range("AC1").Copy($ws.Range("AD1:AF1"))
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row on this sheet belongs to the same team/season, so the
# team's season record (76 wins, 86 losses, 0 ties) is repeated for each
# data row, 2 through 53.
$wins = $ws.Range("AD2:AD53")
$losses = $ws.Range("AE2:AE53")
$ties = $ws.Range("AF2:AF53")

$wins.Value = 76
$losses.Value = 86
$ties.Value = 0
